# Reorders the "Recorded By" (column G) list in each data row so that any
# "System"/"system" entries come first (keeping their relative order), and
# the rest of the names/emails are sorted alphabetically.
#
# Example: "backup@backdoor.com, System, system" -> "System, system, backup@backdoor.com"
#          "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
#          "dnasr281@gmail.com, admin@admin.com"  -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row() + $usedRange.Rows.Count - 1

for ($row = 2; $row -le $lastRow; $row++) {

    $cell = $ws.Cells.Item($row, 7)
    $raw = $cell.Value()

    if ($raw -eq $null) {
        continue
    }

    $text = [string]$raw
    if ($text -eq "") {
        continue
    }

    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $parts = $text.Split(",")

    $systemItems = @()
    $otherItems = @()

    foreach ($part in $parts) {
        $trimmed = $part.Trim()
        if ($trimmed.ToLower() -eq "system") {
            $systemItems += $trimmed
        } else {
            $otherItems += $trimmed
        }
    }

    $otherItemsSorted = $otherItems | Sort-Object

    $result = @()
    $result += $systemItems
    $result += $otherItemsSorted

    $newText = $result -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
